# Regenerate Report for Handback: bump timestamps / priority for the
# c05034c4-8eb9-4f16-8337-e06c142a01c8 file's latest handback cycle.
#
# Because several rows in this fixture happen to share identical text
# (e.g. row 4 and row 5 previously had the same generated timestamp),
# every cell that currently displays the OLD text must be updated so the
# shared-string pool ends up with the same set of distinct values as the
# target workbook.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ---------------------------------------------------
# Column G = "Latest HO Xliff Generate Date"
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G4").Value = "2016-08-28 14:16:56"
$wsOverview.Range("G5").Value = "2016-08-28 14:16:56"

# ---- zh-cn sheet --------------------------------------------------------
# Column E = "Priority", H = "Correspond Handoff Datetime",
# K = "Correspond Handback DateTime"
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E4").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
$wsZhCn.Range("H4").Value = "2016-08-28 14:16:51"
$wsZhCn.Range("H5").Value = "2016-08-28 14:16:51"
$wsZhCn.Range("K4").Value = "2016-08-28 14:17:12"
$wsZhCn.Range("K5").Value = "2016-08-28 14:17:12"

# ---- de-de sheet --------------------------------------------------------
# Column E = "Priority", H = "Correspond Handoff Datetime",
# K = "Correspond Handback DateTime"
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E4").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
$wsDeDe.Range("H4").Value = "2016-08-28 14:16:56"
$wsDeDe.Range("H5").Value = "2016-08-28 14:16:56"
$wsDeDe.Range("K4").Value = "2016-08-28 14:17:19"
$wsDeDe.Range("K5").Value = "2016-08-28 14:17:19"
